$wb = $excel.ActiveWorkbook

# Sheets that share the A:O (panel/measures) layout and get a new column
# inserted immediately before the existing "comment" column (O), pushing
# it out to P and opening up a fresh, blank O column for "item_num".
$layoutSheetNames = @("Measures", "ID", "Dems", "Dates", "NewVars")

foreach ($name in $layoutSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(15).Insert()
    $ws.Range("O1").Value = "item_num"
}

# --- Measures sheet: fill the new item_num column with 1 for every data row ---
$measures = $wb.Worksheets.Item("Measures")
$measures.Range("O2:O185").Value = 1

# Two rows had a stray 30-day time_frame value that didn't belong - clear it.
$measures.Range("K163").ClearContents()
$measures.Range("K164").ClearContents()

# Rows 166:185 carried a leftover font-only style with no visible effect;
# drop it back to the default cell style.
$measures.Range("H166:H185").ClearFormats()

# --- ID / Dems / Dates / NewVars: mirror the divider formatting that sits
# at row 28 on Measures (column N) into the freshly inserted column O ---
$measures.Range("N28").Copy() | Out-Null
$dividerTargets = @("ID", "Dems", "Dates", "NewVars")
foreach ($name in $dividerTargets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("O28").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Workbook-level: the filter database on Measures now spans through P ---
$filterName = $wb.Names.Item("Measures!_FilterDatabase")
$filterName.RefersTo = "=Measures!`$A`$1:`$P`$185"

# --- Selections left by the editing session ---
foreach ($name in $dividerTargets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $ws.Range("O1:O1048576").Select() | Out-Null
}

$measures.Activate()
$measures.Range("N186").Select() | Out-Null
